$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value = 22.51000000000008
$ws.Cells.Item(2, 8).Value = [double]"1.887379141862766e-15"
$ws.Cells.Item(2, 9).Value = [double]"1.887379141862766e-15"
$ws.Cells.Item(2, 12).Value = 47.48749788286049
$ws.Cells.Item(2, 13).Value = "[39.42585729632835, 55.549138469392624]"
$ws.Cells.Item(2, 14).Value = [double]"1.77635683940025e-15"
$ws.Cells.Item(2, 15).Value = [double]"1.77635683940025e-15"
$ws.Cells.Item(2, 17).Value = "[1.3648160276856558, 1.7170266154755023]"
$ws.Cells.Item(2, 20).Value = 52.86215979700135
$ws.Cells.Item(2, 21).Value = "[47.63192753109485, 58.092392062907855]"
$ws.Cells.Item(2, 24).Value = 16.98952952952959
$ws.Cells.Item(2, 25).Value = 16.35861861861868
$ws.Cells.Item(2, 26).Value = 17.6204404404405
$ws.Cells.Item(3, 6).Value = 22.51000000000008
$ws.Cells.Item(3, 12).Value = 45.80650202319914
$ws.Cells.Item(3, 13).Value = "[36.88612207060462, 54.72688197579367]"
$ws.Cells.Item(3, 14).Value = [double]"1.791899961745003e-13"
$ws.Cells.Item(3, 15).Value = [double]"1.791899961745003e-13"
$ws.Cells.Item(3, 16).Value = 2.106974051957119
$ws.Cells.Item(3, 17).Value = "[1.9057108589343503, 2.308237244979887]"
$ws.Cells.Item(3, 20).Value = 51.641537680898
$ws.Cells.Item(3, 21).Value = "[46.82658647366479, 56.456488888131204]"
$ws.Cells.Item(3, 24).Value = 14.96160160160165
$ws.Cells.Item(3, 25).Value = 14.24056056056061
$ws.Cells.Item(3, 26).Value = 15.68264264264269
$ws.Cells.Item(4, 6).Value = 22.51000000000008
$ws.Cells.Item(4, 8).Value = [double]"5.10702591327572e-15"
$ws.Cells.Item(4, 9).Value = [double]"5.10702591327572e-15"
$ws.Cells.Item(4, 12).Value = 38.93573059781723
$ws.Cells.Item(4, 13).Value = "[29.982394463210667, 47.8890667324238]"
$ws.Cells.Item(4, 14).Value = [double]"2.784705799285803e-11"
$ws.Cells.Item(4, 15).Value = [double]"2.784705799285803e-11"
$ws.Cells.Item(4, 16).Value = 2.52207938756658
$ws.Cells.Item(4, 17).Value = "[2.295658295415965, 2.7485004797171957]"
$ws.Cells.Item(4, 20).Value = 49.90611012718713
$ws.Cells.Item(4, 21).Value = "[45.285263074555445, 54.52695717981882]"
$ws.Cells.Item(4, 24).Value = 13.4744544544545
$ws.Cells.Item(4, 25).Value = 12.66328328328333
$ws.Cells.Item(4, 26).Value = 14.28562562562568
$ws.Cells.Item(5, 6).Value = 22.51000000000008
$ws.Cells.Item(5, 12).Value = 50.54645910470322
$ws.Cells.Item(5, 13).Value = "[41.211000942515405, 59.88191726689103]"
$ws.Cells.Item(5, 14).Value = [double]"3.219646771412954e-14"
$ws.Cells.Item(5, 15).Value = [double]"3.219646771412954e-14"
$ws.Cells.Item(5, 16).Value = 2.836553126664658
$ws.Cells.Item(5, 17).Value = "[2.647868883205812, 3.0252373701235036]"
$ws.Cells.Item(5, 20).Value = 51.24363000053734
$ws.Cells.Item(5, 21).Value = "[45.82672173999086, 56.660538261083815]"
$ws.Cells.Item(5, 24).Value = 12.34782782782787
$ws.Cells.Item(5, 25).Value = 11.6718518518519
$ws.Cells.Item(5, 26).Value = 13.02380380380385
$ws.Cells.Item(6, 6).Value = 22.51000000000008
$ws.Cells.Item(6, 8).Value = [double]"1.110223024625157e-16"
$ws.Cells.Item(6, 9).Value = [double]"1.110223024625157e-16"
$ws.Cells.Item(6, 12).Value = 49.60806844949819
$ws.Cells.Item(6, 13).Value = "[43.100016139943925, 56.11612075905246]"
$ws.Cells.Item(6, 14).Value = 0
$ws.Cells.Item(6, 15).Value = 0
$ws.Cells.Item(6, 16).Value = -3.031526844905466
$ws.Cells.Item(6, 17).Value = "[-3.16989529010862, -2.893158399702312]"
$ws.Cells.Item(6, 20).Value = 50.9057633912772
$ws.Cells.Item(6, 21).Value = "[46.700297487863125, 55.11122929469127]"
$ws.Cells.Item(6, 24).Value = 10.86068068068072
$ws.Cells.Item(6, 25).Value = 10.364964964965
$ws.Cells.Item(6, 26).Value = 11.35639639639644
$ws.Cells.Item(7, 6).Value = 22.51000000000008
$ws.Cells.Item(7, 12).Value = 46.59867815397382
$ws.Cells.Item(7, 13).Value = "[37.41029582308299, 55.787060484864654]"
$ws.Cells.Item(7, 14).Value = [double]"2.664535259100376e-13"
$ws.Cells.Item(7, 15).Value = [double]"2.664535259100376e-13"
$ws.Cells.Item(7, 16).Value = -2.591263610168158
$ws.Cells.Item(7, 17).Value = "[-2.792526803190928, -2.3900004171453877]"
$ws.Cells.Item(7, 20).Value = 51.52235779576223
$ws.Cells.Item(7, 21).Value = "[46.79118301853575, 56.25353257298872]"
$ws.Cells.Item(7, 24).Value = 9.283403403403437
$ws.Cells.Item(7, 25).Value = 8.562362362362391
$ws.Cells.Item(7, 26).Value = 10.00444444444448
$ws.Cells.Item(8, 6).Value = 22.40000000000006
$ws.Cells.Item(8, 8).Value = [double]"1.554312234475219e-15"
$ws.Cells.Item(8, 9).Value = [double]"1.554312234475219e-15"
$ws.Cells.Item(8, 12).Value = 40.63923172720357
$ws.Cells.Item(8, 13).Value = "[31.894930934782153, 49.38353251962498]"
$ws.Cells.Item(8, 14).Value = [double]"3.952393967665557e-12"
$ws.Cells.Item(8, 15).Value = [double]"3.952393967665557e-12"
$ws.Cells.Item(8, 16).Value = -2.465474114528927
$ws.Cells.Item(8, 17).Value = "[-2.679316257115619, -2.2516319719422344]"
$ws.Cells.Item(8, 20).Value = 52.91474933140149
$ws.Cells.Item(8, 21).Value = "[48.346377475310256, 57.48312118749272]"
$ws.Cells.Item(8, 24).Value = 8.789589589589616
$ws.Cells.Item(8, 25).Value = 8.027227227227252
$ws.Cells.Item(8, 26).Value = 9.551951951951979
$ws.Cells.Item(9, 6).Value = 22.40000000000006
$ws.Cells.Item(9, 8).Value = [double]"1.076916333886402e-14"
$ws.Cells.Item(9, 9).Value = [double]"1.076916333886402e-14"
$ws.Cells.Item(9, 12).Value = 44.9022848468385
$ws.Cells.Item(9, 13).Value = "[35.24040983342461, 54.5641598602524]"
$ws.Cells.Item(9, 14).Value = [double]"3.956168725949283e-12"
$ws.Cells.Item(9, 15).Value = [double]"3.956168725949283e-12"
$ws.Cells.Item(9, 16).Value = -2.037789829355541
$ws.Cells.Item(9, 17).Value = "[-2.264210921506157, -1.8113687372049245]"
$ws.Cells.Item(9, 20).Value = 53.9394749416461
$ws.Cells.Item(9, 21).Value = "[48.71582168927402, 59.16312819401817]"
$ws.Cells.Item(9, 24).Value = 7.264864864864881
$ws.Cells.Item(9, 25).Value = 6.45765765765767
$ws.Cells.Item(9, 26).Value = 8.072072072072093
$ws.Cells.Item(10, 6).Value = 22.40000000000006
$ws.Cells.Item(10, 8).Value = [double]"2.120525977034049e-14"
$ws.Cells.Item(10, 9).Value = [double]"2.120525977034049e-14"
$ws.Cells.Item(10, 12).Value = 44.34124310131111
$ws.Cells.Item(10, 13).Value = "[36.051141495167954, 52.63134470745426]"
$ws.Cells.Item(10, 14).Value = [double]"4.796163466380676e-14"
$ws.Cells.Item(10, 15).Value = [double]"4.796163466380676e-14"
$ws.Cells.Item(10, 16).Value = -1.610105544182156
$ws.Cells.Item(10, 17).Value = "[-1.7987897876410024, -1.4214213007233099]"
$ws.Cells.Item(10, 20).Value = 48.33837295510011
$ws.Cells.Item(10, 21).Value = "[43.10978775055437, 53.56695815964584]"
$ws.Cells.Item(10, 24).Value = 5.740140140140156
$ws.Cells.Item(10, 25).Value = 5.067467467467481
$ws.Cells.Item(10, 26).Value = 6.41281281281283
$ws.Cells.Item(11, 6).Value = 22.40000000000006
$ws.Cells.Item(11, 8).Value = [double]"2.703393064962256e-13"
$ws.Cells.Item(11, 9).Value = [double]"2.703393064962256e-13"
$ws.Cells.Item(11, 12).Value = 43.75974713979602
$ws.Cells.Item(11, 13).Value = "[33.44493519157367, 54.07455908801837]"
$ws.Cells.Item(11, 14).Value = [double]"5.631894950397509e-11"
$ws.Cells.Item(11, 15).Value = [double]"5.631894950397509e-11"
$ws.Cells.Item(11, 16).Value = -1.144684410317002
$ws.Cells.Item(11, 17).Value = "[-1.3962634015954638, -0.8931054190385392]"
$ws.Cells.Item(11, 18).Value = [double]"7.441158800247649e-12"
$ws.Cells.Item(11, 19).Value = [double]"7.441158800247649e-12"
$ws.Cells.Item(11, 20).Value = 52.19056211821905
$ws.Cells.Item(11, 21).Value = "[46.5110873196461, 57.870036916791996]"
$ws.Cells.Item(11, 24).Value = 4.080880880880891
$ws.Cells.Item(11, 25).Value = 3.18398398398399
$ws.Cells.Item(11, 26).Value = 4.977777777777791
$ws.Cells.Item(12, 6).Value = 22.40000000000006
$ws.Cells.Item(12, 12).Value = 50.68013465967335
$ws.Cells.Item(12, 13).Value = "[40.473072423500774, 60.887196895845925]"
$ws.Cells.Item(12, 14).Value = [double]"5.193623309196482e-13"
$ws.Cells.Item(12, 15).Value = [double]"5.193623309196482e-13"
$ws.Cells.Item(12, 16).Value = -0.6918422260157699
$ws.Cells.Item(12, 17).Value = "[-0.8931054190385392, -0.49057903299300065]"
$ws.Cells.Item(12, 18).Value = [double]"1.323369414052422e-08"
$ws.Cells.Item(12, 19).Value = [double]"1.323369414052422e-08"
$ws.Cells.Item(12, 20).Value = 52.45047808389807
$ws.Cells.Item(12, 21).Value = "[47.14462921034393, 57.75632695745221]"
$ws.Cells.Item(12, 24).Value = 2.466466466466471
$ws.Cells.Item(12, 25).Value = 1.748948948948953
$ws.Cells.Item(12, 26).Value = 3.18398398398399
$ws.Cells.Item(13, 2).Value = 0
$ws.Cells.Item(13, 6).Value = 22.40000000000006
$ws.Cells.Item(13, 8).Value = [double]"1.77635683940025e-15"
$ws.Cells.Item(13, 9).Value = [double]"1.77635683940025e-15"
$ws.Cells.Item(13, 12).Value = 50.28801761002323
$ws.Cells.Item(13, 13).Value = "[41.21332327125748, 59.362711948788984]"
$ws.Cells.Item(13, 14).Value = [double]"1.48769885299771e-14"
$ws.Cells.Item(13, 15).Value = [double]"1.48769885299771e-14"
$ws.Cells.Item(13, 16).Value = -0.1006315965113851
$ws.Cells.Item(13, 17).Value = "[-0.3144737390980774, 0.11321054607530723]"
$ws.Cells.Item(13, 18).Value = 0.3482868319509413
$ws.Cells.Item(13, 19).Value = 0.3482868319509413
$ws.Cells.Item(13, 20).Value = 58.11436452429129
$ws.Cells.Item(13, 21).Value = "[52.24419403129664, 63.98453501728594]"
$ws.Cells.Item(13, 24).Value = 0.3587587587587571
$ws.Cells.Item(13, 25).Value = -0.4036036036036064
$ws.Cells.Item(13, 26).Value = 1.121121121121121
$ws.Cells.Item(14, 6).Value = 22.40000000000006
$ws.Cells.Item(14, 8).Value = [double]"2.015054789694659e-13"
$ws.Cells.Item(14, 9).Value = [double]"2.015054789694659e-13"
$ws.Cells.Item(14, 12).Value = 42.93045040438788
$ws.Cells.Item(14, 13).Value = "[33.16394361741014, 52.69695719136561]"
$ws.Cells.Item(14, 14).Value = [double]"2.043631930348511e-11"
$ws.Cells.Item(14, 15).Value = [double]"2.043631930348511e-11"
$ws.Cells.Item(14, 16).Value = 0.2830263651882703
$ws.Cells.Item(14, 17).Value = "[0.01886842434588587, 0.5471843060306547]"
$ws.Cells.Item(14, 18).Value = 0.0363044836063966
$ws.Cells.Item(14, 19).Value = 0.0363044836063966
$ws.Cells.Item(14, 20).Value = 53.76429199794611
$ws.Cells.Item(14, 21).Value = "[48.07005328790434, 59.45853070798788]"
$ws.Cells.Item(14, 24).Value = 21.39099099099105
$ws.Cells.Item(14, 25).Value = 20.44924924924931
$ws.Cells.Item(14, 26).Value = 22.33273273273279
$ws.Cells.Item(15, 6).Value = 22.40000000000006
$ws.Cells.Item(15, 8).Value = [double]"1.110223024625157e-16"
$ws.Cells.Item(15, 9).Value = [double]"1.110223024625157e-16"
$ws.Cells.Item(15, 12).Value = 49.16599020767087
$ws.Cells.Item(15, 13).Value = "[40.816654784691735, 57.51532563065]"
$ws.Cells.Item(15, 14).Value = [double]"1.77635683940025e-15"
$ws.Cells.Item(15, 15).Value = [double]"1.77635683940025e-15"
$ws.Cells.Item(15, 16).Value = 0.5597632555945786
$ws.Cells.Item(15, 17).Value = "[0.38365796169965627, 0.735868549489501]"
$ws.Cells.Item(15, 18).Value = [double]"7.887318398402954e-08"
$ws.Cells.Item(15, 19).Value = [double]"7.887318398402954e-08"
$ws.Cells.Item(15, 20).Value = 50.68690935667668
$ws.Cells.Item(15, 21).Value = "[46.13626127611355, 55.23755743723981]"
$ws.Cells.Item(15, 24).Value = 20.40440440440446
$ws.Cells.Item(15, 25).Value = 19.77657657657663
$ws.Cells.Item(15, 26).Value = 21.03223223223229
